$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the cells we are about to update so that
# numeric-looking strings (percentages, plain numbers, decimals) are
# preserved verbatim as text, matching the original inline-string cells.
$cellRefs = @(
    "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7",
    "G7", "D8", "E8", "G8", "B9", "C9", "D9", "E9", "G9", "B10", "C10", "D10", "E10", "G10", "B11",
    "C11", "D11", "E11", "G11", "B12", "C12", "D12", "E12", "G12", "B13", "C13", "D13", "E13", "G13",
    "B14", "C14", "D14", "E14", "G14", "B15", "C15", "D15", "E15", "G15", "B16", "C16", "D16", "E16",
    "G16", "B17", "C17", "D17", "E17", "G17", "B18", "C18", "D18", "E18", "G18", "B19", "C19", "D19",
    "E19", "G19", "B20", "C20", "D20", "E20", "G20", "D21", "E21", "G21", "E22", "G22", "D23", "E23",
    "G23", "E24", "G24", "E25", "G25", "D26", "E26", "G26", "E27", "G27", "D28", "E28", "G28", "G29",
    "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "G39", "D40", "E40", "G40", "B41",
    "C41", "D41", "E41", "G41", "B42", "C42", "D42", "E42", "G42", "E43", "G43", "D44", "E44", "G44",
    "D45", "E45", "G45", "E46", "G46", "E47", "G47", "D48", "E48", "G48", "E49", "G49", "E50", "G50",
    "G51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("E2").Value = '-1.09%'
$ws.Range("G2").Value = '3'
$ws.Range("D3").Value = '29.13'
$ws.Range("E3").Value = '9.85%'
$ws.Range("G3").Value = '3'
$ws.Range("D4").Value = '5.107'
$ws.Range("E4").Value = '-0.48%'
$ws.Range("G4").Value = '3'
$ws.Range("D5").Value = '0.05649'
$ws.Range("E5").Value = '1.14%'
$ws.Range("G5").Value = '3'
$ws.Range("D6").Value = '6.495'
$ws.Range("E6").Value = '0.31%'
$ws.Range("G6").Value = '3'
$ws.Range("D7").Value = '0.8268'
$ws.Range("E7").Value = '1.24%'
$ws.Range("G7").Value = '3'
$ws.Range("D8").Value = '0.8660'
$ws.Range("E8").Value = '3.01%'
$ws.Range("G8").Value = '3'
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").Value = '0.1329'
$ws.Range("E9").Value = '-0.02%'
$ws.Range("G9").Value = '3'
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D10").Value = '0.02859'
$ws.Range("E10").Value = '-0.24%'
$ws.Range("G10").Value = '3'
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D11").Value = '0.09374'
$ws.Range("E11").Value = '-0.20%'
$ws.Range("G11").Value = '3'
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D12").Value = '0.001509'
$ws.Range("E12").Value = '-0.99%'
$ws.Range("G12").Value = '3'
$ws.Range("B13").Value = 'CoinExToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D13").Value = '0.04145'
$ws.Range("E13").Value = '-9.79%'
$ws.Range("G13").Value = '3'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").Value = '0.0006018'
$ws.Range("E14").Value = '0.65%'
$ws.Range("G14").Value = '3'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006115'
$ws.Range("E15").Value = '-1.12%'
$ws.Range("G15").Value = '3'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.521'
$ws.Range("E16").Value = '-3.16%'
$ws.Range("G16").Value = '3'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '3.021'
$ws.Range("E17").Value = '-0.61%'
$ws.Range("G17").Value = '3'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").Value = '2.219'
$ws.Range("E18").Value = '1.66%'
$ws.Range("G18").Value = '3'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").Value = '0.3149'
$ws.Range("E19").Value = '1.19%'
$ws.Range("G19").Value = '3'
$ws.Range("B20").Value = 'MandalaExchangeToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D20").Value = '0.06929'
$ws.Range("E20").Value = '-0.99%'
$ws.Range("G20").Value = '3'
$ws.Range("D21").Value = '0.03245'
$ws.Range("E21").Value = '5.83%'
$ws.Range("G21").Value = '3'
$ws.Range("E22").Value = '-2.03%'
$ws.Range("G22").Value = '3'
$ws.Range("D23").Value = '3.632'
$ws.Range("E23").Value = '-3.24%'
$ws.Range("G23").Value = '3'
$ws.Range("E24").Value = '-0.01%'
$ws.Range("G24").Value = '3'
$ws.Range("E25").Value = '-2.99%'
$ws.Range("G25").Value = '3'
$ws.Range("D26").Value = '0.004443'
$ws.Range("E26").Value = '-1.79%'
$ws.Range("G26").Value = '3'
$ws.Range("E27").Value = '22.79%'
$ws.Range("G27").Value = '3'
$ws.Range("D28").Value = '0.0001404'
$ws.Range("E28").Value = '0.56%'
$ws.Range("G28").Value = '3'
$ws.Range("G29").Value = '3'
$ws.Range("G30").Value = '3'
$ws.Range("G31").Value = '3'
$ws.Range("G32").Value = '3'
$ws.Range("G33").Value = '3'
$ws.Range("G34").Value = '3'
$ws.Range("G35").Value = '3'
$ws.Range("G36").Value = '3'
$ws.Range("G37").Value = '3'
$ws.Range("G38").Value = '3'
$ws.Range("G39").Value = '3'
$ws.Range("D40").Value = '0.03708'
$ws.Range("E40").Value = '1.77%'
$ws.Range("G40").Value = '3'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1054'
$ws.Range("E41").Value = '-23.69%'
$ws.Range("G41").Value = '3'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '0.003416'
$ws.Range("E42").Value = '-44.49%'
$ws.Range("G42").Value = '3'
$ws.Range("E43").Value = '-12.18%'
$ws.Range("G43").Value = '3'
$ws.Range("D44").Value = '0.009464'
$ws.Range("E44").Value = '5.62%'
$ws.Range("G44").Value = '3'
$ws.Range("D45").Value = '0.00005105'
$ws.Range("E45").Value = '-4.35%'
$ws.Range("G45").Value = '3'
$ws.Range("E46").Value = '-0.01%'
$ws.Range("G46").Value = '3'
$ws.Range("E47").Value = '-3.68%'
$ws.Range("G47").Value = '3'
$ws.Range("D48").Value = '0.002506'
$ws.Range("E48").Value = '-1.64%'
$ws.Range("G48").Value = '3'
$ws.Range("E49").Value = '-0.01%'
$ws.Range("G49").Value = '3'
$ws.Range("E50").Value = '-0.01%'
$ws.Range("G50").Value = '3'
$ws.Range("G51").Value = '3'

# Restore default (Normal) style so no stray number-format styling remains
foreach ($ref in $cellRefs) {
    $ws.Range($ref).Style = "Normal"
}
